$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.340.09'
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('D3').Value = '1.592.55'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '210.25'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.507'
$ws.Range('E6').Value = '  -1.14%  '
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.60'
$ws.Range('E10').Value = '  -0.52%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').Value = '1.814.30'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.08'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.562.70'
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.519'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('D17').Value = '26.326.84'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.50'
$ws.Range('E19').Value = '  +4.28%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '211.94'
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.19'
$ws.Range('E23').Value = '  -1.99%  '
$ws.Range('E24').Value = '  -1.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.41'
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('E27').Value = '  -1.33%  '
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.28'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('E30').Value = '  -0.91%  '
$ws.Range('E31').Value = '  -0.41%  '
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('D34').Value = '1.302.19'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.614'
$ws.Range('E35').Value = '  +3.55%  '
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.10'
$ws.Range('E39').Value = '  -14.38%  '
$ws.Range('E40').Value = '  -1.81%  '
$ws.Range('E41').Value = '  -0.51%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.63'
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '62.71'
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.13'
$ws.Range('E44').Value = '  -2.58%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.761'
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('D46').Value = '1.728.00'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '88.51'
$ws.Range('E47').Value = '  -2.20%  '
$ws.Range('E48').Value = '  -3.78%  '
$ws.Range('D49').Value = '0.0₆0103'
$ws.Range('E49').Value = '  -2.48%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0985'
$ws.Range('E50').Value = '  -3.97%  '
$ws.Range('E51').Value = '  -1.42%  '
